# DeveloperGuide: remove child Command classes in class diagram
#
# On the "Diagrams" slide (slide 2) the class-diagram for the Command
# hierarchy showed the "Command" box fanning out (via a small triangle
# junction + elbow connectors) to its concrete children: AddCommand,
# ClearCommand ("CrearCommand"), IncorrectCommand, and an "...Command"
# box standing in for "more commands". This change deletes that whole
# fan-out - the three named child boxes, the "...Command" placeholder
# box, the junction triangle, and all four elbow connectors wiring them
# to the triangle - leaving just the "Command" box itself.
#
# Shapes removed (by their p:cNvPr id / name):
#   50 "Rectangle 49"          -> AddCommand
#   52 "Rectangle 51"          -> CrearCommand
#   53 "Rectangle 52"          -> IncorrectCommand
#   54 "Isosceles Triangle 53" -> fan-out junction triangle
#   55 "Elbow Connector 54"    -> triangle -> CrearCommand
#   56 "Elbow Connector 55"    -> triangle -> AddCommand
#   57 "Elbow Connector 56"    -> triangle -> IncorrectCommand
#   59 "Rectangle 58"          -> "...Command"
#   60 "Elbow Connector 59"    -> triangle -> "...Command"
#
# The "Command" box (id 46) and the dashed connector from the UniqueTagList
# rectangle into "Command" (id 67, "Elbow Connector 66") are untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$idsToRemove = @(50, 52, 53, 54, 55, 56, 57, 59, 60)

foreach ($id in $idsToRemove) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shape = $s.Shapes.Item($i)
        if ($shape.Id -eq $id) {
            $shape.Delete()
            break
        }
    }
}
